$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'275.62"
$ws.Range("E2").Value = "'-1.04%"
$ws.Range("D3").Value = "'27.31"
$ws.Range("E3").Value = "'1.30%"
$ws.Range("D4").Value = "'4.788"
$ws.Range("E4").Value = "'-3.05%"
$ws.Range("E5").Value = "'-1.09%"
$ws.Range("D6").Value = "'6.942"
$ws.Range("E6").Value = "'-0.95%"
$ws.Range("D7").Value = "'1.347"
$ws.Range("E7").Value = "'29.77%"
$ws.Range("D8").Value = "'0.8773"
$ws.Range("E8").Value = "'-1.21%"
$ws.Range("D9").Value = "'0.1513"
$ws.Range("E9").Value = "'1.34%"
$ws.Range("D10").Value = "'0.05080"
$ws.Range("E10").Value = "'-1.90%"
$ws.Range("D11").Value = "'0.07583"
$ws.Range("E11").Value = "'2.59%"
$ws.Range("D12").Value = "'0.02963"
$ws.Range("E12").Value = "'-4.17%"
$ws.Range("D13").Value = "'0.09025"
$ws.Range("E13").Value = "'-0.56%"
$ws.Range("D14").Value = "'0.001578"
$ws.Range("E14").Value = "'1.09%"
$ws.Range("D15").Value = "'0.0006386"
$ws.Range("E15").Value = "'0.35%"
$ws.Range("D16").Value = "'0.005903"
$ws.Range("E16").Value = "'-1.89%"
$ws.Range("D17").Value = "'3.450"
$ws.Range("E17").Value = "'-1.37%"
$ws.Range("E18").Value = "'-1.69%"
$ws.Range("E19").Value = "'-1.04%"
$ws.Range("E20").Value = "'-1.21%"
$ws.Range("D21").Value = "'0.1344"
$ws.Range("E21").Value = "'0.80%"
$ws.Range("D22").Value = "'3.904"
$ws.Range("E22").Value = "'-0.55%"
$ws.Range("D23").Value = "'0.04399"
$ws.Range("E23").Value = "'0.81%"
$ws.Range("D24").Value = "'0.001170"
$ws.Range("E24").Value = "'-1.27%"
$ws.Range("D25").Value = "'0.003862"
$ws.Range("E25").Value = "'4.66%"
$ws.Range("E26").Value = "'-0.42%"
$ws.Range("E27").Value = "'13.89%"
$ws.Range("D40").Value = "'0.04160"
$ws.Range("E40").Value = "'1.82%"
$ws.Range("D41").Value = "'0.006849"
$ws.Range("E41").Value = "'2.55%"
$ws.Range("D42").Value = "'0.1177"
$ws.Range("E42").Value = "'-0.07%"
$ws.Range("D43").Value = "'0.002038"
$ws.Range("E43").Value = "'-13.92%"
$ws.Range("D44").Value = "'0.01153"
$ws.Range("E44").Value = "'-8.12%"
$ws.Range("D45").Value = "'0.00005165"
$ws.Range("D46").Value = "'1.486"
$ws.Range("E46").Value = "'-36.89%"
$ws.Range("E47").Value = "'2.31%"
